$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "department" column (C) values for the data rows.
# Rows 2-7: individual courses -> "Automotive"
# Rows 8-9: course packages -> "Packages"
$ws.Range("C2:C7").Value = "Automotive"
$ws.Range("C8:C9").Value = "Packages"
